$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INCO")

# The old row 35 (Caso -509) was removed; every subsequent row shifted up by one.
$ws.Rows.Item(35).Delete()
